$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.762.44"
$ws.Range("E2").Value = "  +4.33%  "
$ws.Range("D3").Value = "'2.421.99"
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'316.46"
$ws.Range("E5").Value = "  +4.51%  "
$ws.Range("D6").Value = "'101.49"
$ws.Range("E6").Value = "  +6.68%  "
$ws.Range("E7").Value = "  +2.29%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.530"
$ws.Range("E9").Value = "  +11.49%  "
$ws.Range("D10").Value = "'35.38"
$ws.Range("E10").Value = "  +3.28%  "
$ws.Range("D11").Value = "'0.0799"
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").Value = "'18.75"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("E14").Value = "  +3.08%  "
$ws.Range("D15").Value = "'2.799.41"
$ws.Range("E15").Value = "  +2.72%  "
$ws.Range("D16").Value = "'2.417.61"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "'0.832"
$ws.Range("E17").Value = "  +4.77%  "
$ws.Range("D18").Value = "'44.551.88"
$ws.Range("E18").Value = "  +3.78%  "
$ws.Range("D19").Value = "'12.32"
$ws.Range("E19").Value = "  +3.13%  "
$ws.Range("D20").Value = "'6.37"
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("D21").Value = "'0.0₃0916"
$ws.Range("E21").Value = "  +3.62%  "
$ws.Range("D22").Value = "'68.74"
$ws.Range("E22").Value = "  +1.25%  "
$ws.Range("D23").Value = "'242.56"
$ws.Range("E23").Value = "  +2.99%  "
$ws.Range("D24").Value = "'2.28"
$ws.Range("E24").Value = "  +4.58%  "
$ws.Range("E25").Value = "  +2.36%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "'25.20"
$ws.Range("E27").Value = "  +3.51%  "
$ws.Range("D28").Value = "'2.25"
$ws.Range("E28").Value = "  -4.99%  "
$ws.Range("D29").Value = "'9.48"
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("D30").Value = "'33.66"
$ws.Range("E30").Value = "  +3.87%  "
$ws.Range("D31").Value = "'48.53"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("E32").Value = "  +18.13%  "
$ws.Range("D33").Value = "'19.52"
$ws.Range("E33").Value = "  +11.63%  "
$ws.Range("D34").Value = "'5.16"
$ws.Range("E34").Value = "  +3.22%  "
$ws.Range("D35").Value = "'0.0772"
$ws.Range("E35").Value = "  +6.61%  "
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("E37").Value = "  +3.67%  "
$ws.Range("D38").Value = "'4.48"
$ws.Range("E38").Value = "  +3.49%  "
$ws.Range("D39").Value = "'2.85"
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("D40").Value = "'123.13"
$ws.Range("E40").Value = "  -4.30%  "
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("D42").Value = "'2.19"
$ws.Range("E42").Value = "  -3.05%  "
$ws.Range("D43").Value = "'21.12"
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("D44").Value = "'0.0291"
$ws.Range("E44").Value = "  +4.46%  "
$ws.Range("D45").Value = "'1.941.48"
$ws.Range("E45").Value = "  +0.78%  "
$ws.Range("E46").Value = "  -1.61%  "
$ws.Range("D47").Value = "'2.94"
$ws.Range("E47").Value = "  +8.11%  "
$ws.Range("D48").Value = "'9.24"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "'1.74"
$ws.Range("E49").Value = "  +15.36%  "
$ws.Range("D50").Value = "'75.76"
$ws.Range("E50").Value = "  +6.13%  "
$ws.Range("D51").Value = "'53.95"
$ws.Range("E51").Value = "  +5.59%  "
